$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("module")

# Connector on module flipped: swap which capacitors go with which pins
$ws.Range("D10").Value = "C1, C2, C5"
$ws.Range("D9").Value = "C6, C7, C8, C9"

$ws.Range("D9").Select()
